# "correct 229v1 growth rates"
# Fixes the "time" column (D3:D13) on the "#2229v1" sheet, which held a
# placeholder 1..11 sequence and should contain the real (irregular)
# elapsed-time values. Also restores the workbook/sheet selection state
# that results from reviewing/selecting that corrected column.

$wb = $excel.ActiveWorkbook

# --- 1. Correct the growth-rate "time" values on sheet "#2229v1" ---------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("D3").Value  = 15.5
$ws1.Range("D4").Value  = 17.25
$ws1.Range("D5").Value  = 19.25
$ws1.Range("D6").Value  = 21.25
$ws1.Range("D7").Value  = 23.25
$ws1.Range("D8").Value  = 25.25
$ws1.Range("D9").Value  = 42.25
$ws1.Range("D10").Value = 45
$ws1.Range("D11").Value = 48
$ws1.Range("D12").Value = 51
$ws1.Range("D13").Value = 66.25

# --- 2. Make "#2229v1" the active sheet/tab, with D2:D13 selected --------
# (the corrected column), matching the reviewer's selection after the fix.
$ws1.Activate()
$ws1.Select()
$ws1.Range("D2:D13").Select()

# --- 3. The other sheets keep their own selections; the previously
# active sheet ("MasterLV3_130v2") is no longer the selected tab, which
# Activate()/Select() above on sheet 1 already took care of.
